$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 577 (shifts existing rows 577-686 down to 579-688)
$ws.Rows("577:578").Insert()

# Row 577 (Primera quality) - new data point
$ws.Range("A577").Value = 8
$ws.Range("B577").Value = "Terminal La Palmera de La Serena"
$ws.Range("C577").Value = "Coquimbo"
$ws.Range("D577").Value = 44694
$ws.Range("E577").Value = 4
$ws.Range("F577").Value = 100112008
$ws.Range("G577").Value = "Coliflor"
$ws.Range("H577").Value = "Sin especificar"
$ws.Range("I577").Value = "Primera"
$ws.Range("J577").Value = 2520
$ws.Range("K577").Value = 850
$ws.Range("L577").Value = 900
$ws.Range("M577").Value = 875
$ws.Range("N577").Value = "$/unidad"
$ws.Range("O577").Value = "Provincia del Elquí"
$ws.Range("P577").Value = 875
$ws.Range("Q577").Value = 1
$ws.Range("R577").Value = "Hortaliza"

# Row 578 (Segunda quality) - new data point
$ws.Range("A578").Value = 8
$ws.Range("B578").Value = "Terminal La Palmera de La Serena"
$ws.Range("C578").Value = "Coquimbo"
$ws.Range("D578").Value = 44694
$ws.Range("E578").Value = 4
$ws.Range("F578").Value = 100112008
$ws.Range("G578").Value = "Coliflor"
$ws.Range("H578").Value = "Sin especificar"
$ws.Range("I578").Value = "Segunda"
$ws.Range("J578").Value = 1320
$ws.Range("K578").Value = 750
$ws.Range("L578").Value = 800
$ws.Range("M578").Value = 775
$ws.Range("N578").Value = "$/unidad"
$ws.Range("O578").Value = "Provincia del Elquí"
$ws.Range("P578").Value = 775
$ws.Range("Q578").Value = 1
$ws.Range("R578").Value = "Hortaliza"
